# Inserts a new data row at row 385 (pushing the existing rows 385-473 down
# to 386-474) and populates it with a new price record for Acelga.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(385).Insert()

$ws.Range("A385").Value = 3
$ws.Range("B385").Value = "Femacal de La Calera"
$ws.Range("C385").Value = "Coquimbo"
$ws.Range("D385").Value = 44943
$ws.Range("E385").Value = 5
$ws.Range("F385").Value = 100112009
$ws.Range("G385").Value = "Acelga"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 160
$ws.Range("K385").Value = 4000
$ws.Range("L385").Value = 4300
$ws.Range("M385").Value = 4094
$ws.Range("N385").Value = "$/docena de atados (6 kilos)"
$ws.Range("O385").Value = "Provincia de Quillota"
$ws.Range("P385").Value = 682
$ws.Range("Q385").Value = 6
$ws.Range("R385").Value = "Hortaliza"
